$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        Row = 185
        Title = "Senior Golang Developer"
        Url = "https://www.dice.com/job-detail/9c7c1315-7b19-4b3c-b881-40d93c6072c9?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
        Location = "Plano, Texas"
        EmploymentType = "Contract"
        Salary = "Depends on Experience"
        Company = "NimbusAITech LLC"
    },
    @{
        Row = 186
        Title = "Golang Architect / Principal Backend Architect Only Local to GA"
        Url = "https://www.dice.com/job-detail/4d3b58b3-db40-4c45-8c13-5474c420def8?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
        Location = "Atlanta, Georgia"
        EmploymentType = "Third Party"
        Salary = "USD85 - USD110"
        Company = "Randstad Digital"
    },
    @{
        Row = 187
        Title = "React Developer (React, Node, AWS; Golang a plus)"
        Url = "https://www.dice.com/job-detail/f427a600-410d-442e-ae37-e687409eabcc?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
        Location = "McLean, Virginia"
        EmploymentType = "Contract"
        Salary = "55 - 60"
        Company = "Stellar Professionals LLC"
    }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Title
    $ws.Cells.Item($r.Row, 2).Value = $r.Url
    $ws.Cells.Item($r.Row, 3).Value = $r.Location
    $ws.Cells.Item($r.Row, 4).Value = $r.EmploymentType
    $ws.Cells.Item($r.Row, 5).Value = $r.Salary
    $ws.Cells.Item($r.Row, 6).Value = $r.Company
}
